$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "33.954.89"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3
$ws.Range("D3").Value = "1.779.35"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.08"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.22"
$ws.Range("E8").Value = "  +2.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.290"
$ws.Range("E9").Value = "  +0.90%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0704"
$ws.Range("E10").Value = "  -0.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  +1.42%  "

# Row 12
$ws.Range("D12").Value = "2.036.31"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("D13").Value = "1.787.77"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.94"
$ws.Range("E14").Value = "  +3.09%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.620"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.943.06"
$ws.Range("E16").Value = "  -0.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.14"
$ws.Range("E17").Value = "  -1.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.86"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.86"
$ws.Range("E19").Value = "  -0.53%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  +1.14%  "

# Row 21
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  +0.54%  "

# Row 23
$ws.Range("E23").Value = "  -0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -2.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.32"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.08"
$ws.Range("E27").Value = "  +0.87%  "

# Row 28
$ws.Range("E28").Value = "  +0.70%  "

# Row 29
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("E30").Value = "  +3.39%  "

# Row 31
$ws.Range("E31").Value = "  -1.13%  "

# Row 32
$ws.Range("E32").Value = "  -0.85%  "

# Row 33
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("E34").Value = "  -0.84%  "

# Row 35
$ws.Range("D35").Value = "1.391.88"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.656"
$ws.Range("E36").Value = "  +4.82%  "

# Row 37
$ws.Range("E37").Value = "  -1.33%  "

# Row 38
$ws.Range("E38").Value = "  +0.85%  "

# Row 39
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.21"
$ws.Range("E40").Value = "  +4.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.910"
$ws.Range("E41").Value = "  -2.01%  "

# Row 42
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "77.63"
$ws.Range("E43").Value = "  -2.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.21"
$ws.Range("E44").Value = "  +13.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.07"
$ws.Range("E45").Value = "  +2.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.11"
$ws.Range("E46").Value = "  +2.22%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0135"
$ws.Range("E47").Value = "  +17.07%  "

# Row 48
$ws.Range("B48").Value = "Kaspa"
$ws.Range("C48").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0496"
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.82"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("D50").Value = "1.936.39"
$ws.Range("E50").Value = "  +0.32%  "

# Row 51
$ws.Range("E51").Value = "  +0.36%  "
